$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("A2").Value = "AIBT Courses Fees 2021.pdf"
$ws.Range("B2").Value = "https://github.com/Viskee-Consultancy/Viskee-Consultancy-Configuration/raw/master/brochures/pdf/onshore/non-coe/aibt/AIBT_Courses_Fees_2021_VOL_2.2.pdf"

# Update row 3
$ws.Range("A3").Value = "AIBT Non-COE Q4 Promotion Brochure.pdf"
$ws.Range("B3").Value = "https://github.com/Viskee-Consultancy/Viskee-Consultancy-Configuration/raw/master/brochures/pdf/onshore/non-coe/aibt/AIBTNon-CoEQ4Brochure_1OCT-31DEC21_VOL1.0.pdf"

# Delete rows 4, 5, 6 entirely (shrink dimension to A1:B3)
$ws.Rows("4:6").Delete()

# Leave the selection where the deleted rows used to be
$ws.Rows("4:6").Select() | Out-Null
